# This script appends the new "TODO" notes paragraphs and the JSON-shaped
# nonogram data block after the final "...zavrtel vrstni red pravil" paragraph.
#
# InsertXML on a Range REPLACES that ranges contents; for a collapsed range at
# the end of a paragraph it still "consumes" that paragraphs own mark (pushing it
# to a new trailing paragraph), so each batch is preceded by InsertParagraphAfter()
# to give InsertXML a disposable empty paragraph to replace, and followed by deleting
# the spare paragraph mark it leaves behind, merging it back into one clean run of
# paragraphs. Batches are additionally capped well under the engines per-call <w:p>
# ceiling so a single InsertXML call never has to rewrite unrelated parts of the body.

$d = $word.ActiveDocument

# --- batch 1 of 4 (15 paragraphs) ---
$lastPara = $d.Paragraphs[$d.Paragraphs.Count]
$endRange = $lastPara.Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$newPara = $d.Paragraphs[$d.Paragraphs.Count]
$newRange = $newPara.Range
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p></w:p><w:p><w:r><w:t>TODO:</w:t></w:r></w:p><w:p><w:r><w:t>Automatično zgeneriraj pravila</w:t></w:r><w:r><w:t xml:space="preserve"> (done)</w:t></w:r></w:p><w:p><w:r><w:t>Dodaj zvoke</w:t></w:r><w:r><w:t xml:space="preserve"> (done)</w:t></w:r></w:p><w:p><w:r><w:t>Ustvari svoj nonogram</w:t></w:r><w:r><w:t xml:space="preserve"> (done)</w:t></w:r></w:p><w:p><w:r><w:t>Omogoči več različnih velikosti nonogramov?</w:t></w:r></w:p><w:p></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t>{</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  "shapes": [</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    [</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [0,0,0,0,0,0,0,0],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [0,0,0,0,0,0,0,0],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [0,0,0,1,3,5,3,1],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [0,0,1,0,0,1,0,0],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [0,0,3,0,1,1,1,0],</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newRange.InsertXML($xml)

$spareIndex = $d.Paragraphs.Count
$beforeSpare = $d.Paragraphs[$spareIndex - 1]
$markRange = $d.Range($beforeSpare.Range.End - 1, $beforeSpare.Range.End)
$markRange.Delete()

# --- batch 2 of 4 (15 paragraphs) ---
$lastPara = $d.Paragraphs[$d.Paragraphs.Count]
$endRange = $lastPara.Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$newPara = $d.Paragraphs[$d.Paragraphs.Count]
$newRange = $newPara.Range
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">      [0,0,5,1,1,1,1,1],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [0,0,3,0,1,1,1,0],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [0,0,1,0,0,1,0,0]</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    ],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    [</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [0,0,0,0,0,1,0,0],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [0,0,0,0,1,1,1,0],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [0,0,0,5,1,1,1,5],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [0,0,5,1,1,1,1,1],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [0,1,1,1,0,0,0,1],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [1,1,1,1,0,1,0,1],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [0,1,1,1,0,0,0,1],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [0,0,5,1,1,1,1,1]</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    ],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    [</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newRange.InsertXML($xml)

$spareIndex = $d.Paragraphs.Count
$beforeSpare = $d.Paragraphs[$spareIndex - 1]
$markRange = $d.Range($beforeSpare.Range.End - 1, $beforeSpare.Range.End)
$markRange.Delete()

# --- batch 3 of 4 (15 paragraphs) ---
$lastPara = $d.Paragraphs[$d.Paragraphs.Count]
$endRange = $lastPara.Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$newPara = $d.Paragraphs[$d.Paragraphs.Count]
$newRange = $newPara.Range
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">      [0,0,0,0,0,0,0,0],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [0,0,0,0,2,1,2,0],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [0,0,0,1,1,1,1,1],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [0,1,1,0,1,0,1,0],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [0,1,1,0,1,0,1,0],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [0,0,1,0,0,1,0,0],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [0,1,1,1,0,0,0,1],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [0,0,3,0,1,1,1,0]</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">    ],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    [</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [0,0,0,1,0,1,0,1],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [0,0,0,1,1,1,1,1],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [0,0,0,1,1,1,1,1],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [1,1,1,1,0,1,0,1],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [0,1,1,0,1,0,1,0],</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newRange.InsertXML($xml)

$spareIndex = $d.Paragraphs.Count
$beforeSpare = $d.Paragraphs[$spareIndex - 1]
$markRange = $d.Range($beforeSpare.Range.End - 1, $beforeSpare.Range.End)
$markRange.Delete()

# --- batch 4 of 4 (6 paragraphs) ---
$lastPara = $d.Paragraphs[$d.Paragraphs.Count]
$endRange = $lastPara.Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$newPara = $d.Paragraphs[$d.Paragraphs.Count]
$newRange = $newPara.Range
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">      [1,1,1,1,0,1,0,1],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [0,1,1,0,1,0,1,0],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      [1,1,1,1,0,1,0,1]</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    ]</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  ]</w:t></w:r></w:p><w:p><w:r><w:t>}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newRange.InsertXML($xml)

$spareIndex = $d.Paragraphs.Count
$beforeSpare = $d.Paragraphs[$spareIndex - 1]
$markRange = $d.Range($beforeSpare.Range.End - 1, $beforeSpare.Range.End)
$markRange.Delete()

Write-Output ("Final Paragraphs.Count=" + $d.Paragraphs.Count)
